$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values stored as text (some look numeric, e.g. "298.51").
# Force text format on the whole Price column first so Excel does not silently
# coerce these into numbers (which would also lose exact formatting/precision),
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '46.773.18'
$ws.Range('D3').Value = '2.308.02'
$ws.Range('D5').Value = '298.51'
$ws.Range('D6').Value = '101.71'
$ws.Range('D10').Value = '36.58'
$ws.Range('D11').Value = '0.0803'
$ws.Range('D12').Value = '7.39'
$ws.Range('D14').Value = '2.655.66'
$ws.Range('D15').Value = '2.305.44'
$ws.Range('D16').Value = '13.97'
$ws.Range('D18').Value = '46.712.73'
$ws.Range('D19').Value = '13.28'
$ws.Range('D22').Value = '67.09'
$ws.Range('D23').Value = '248.26'
$ws.Range('D25').Value = '1.98'
$ws.Range('D26').Value = '1.00'
$ws.Range('D27').Value = '42.82'
$ws.Range('D29').Value = '9.92'
$ws.Range('D31').Value = '5.75'
$ws.Range('D32').Value = '147.03'
$ws.Range('D33').Value = '0.0802'
$ws.Range('D34').Value = '2.62'
$ws.Range('D36').Value = '0.112'
$ws.Range('D39').Value = '15.86'
$ws.Range('D40').Value = '4.04'
$ws.Range('D41').Value = '3.45'
$ws.Range('D43').Value = '2.01'
$ws.Range('D44').Value = '0.999'
$ws.Range('D45').Value = '1.843.30'
$ws.Range('D46').Value = '89.91'
$ws.Range('D48').Value = '75.20'
$ws.Range('D50').Value = '97.44'
$ws.Range('D51').Value = '54.48'

$ws.Range("D2:D51").Style = "Normal"

# Coin name / link / volume columns are plain text already - no coercion risk.
$ws.Range('E2').Value = '  +6.59%  '
$ws.Range('E3').Value = '  +5.16%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('E6').Value = '  +15.02%  '
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('E9').Value = '  +10.22%  '
$ws.Range('E10').Value = '  +13.84%  '
$ws.Range('E11').Value = '  +4.55%  '
$ws.Range('E12').Value = '  +9.36%  '
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('E14').Value = '  +4.99%  '
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('E16').Value = '  +7.53%  '
$ws.Range('E17').Value = '  +6.57%  '
$ws.Range('E18').Value = '  +7.32%  '
$ws.Range('E19').Value = '  +24.29%  '
$ws.Range('E20').Value = '  +6.53%  '
$ws.Range('E21').Value = '  +5.82%  '
$ws.Range('E22').Value = '  +6.75%  '
$ws.Range('E23').Value = '  +7.84%  '
$ws.Range('E24').Value = '  +6.42%  '
$ws.Range('E25').Value = '  +8.66%  '
$ws.Range('E27').Value = '  +19.85%  '
$ws.Range('E28').Value = '  +5.46%  '
$ws.Range('E29').Value = '  +8.22%  '
$ws.Range('E30').Value = '  +5.96%  '
$ws.Range('E31').Value = '  +8.82%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('E33').Value = '  +9.48%  '
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('E35').Value = '  +9.20%  '
$ws.Range('E36').Value = '  +10.49%  '
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range('E38').Value = '  +9.51%  '
$ws.Range('E39').Value = '  +22.73%  '
$ws.Range('E40').Value = '  +15.46%  '
$ws.Range('E41').Value = '  +12.52%  '
$ws.Range('E42').Value = '  +9.50%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E43').Value = '  +24.45%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('E45').Value = '  +3.87%  '
$ws.Range('E46').Value = '  +24.14%  '
$ws.Range('E47').Value = '  +14.89%  '
$ws.Range('E48').Value = '  +16.90%  '
$ws.Range('E49').Value = '  +10.17%  '
$ws.Range('E50').Value = '  +6.48%  '
$ws.Range('E51').Value = '  +11.86%  '
